$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs","Npnt","Itgb1","ECs",3.0,1.0,2.128539,6.385617,0.4682231486221785,0.4682231486221785,3.0,1.0,159.4836373333333,478.450912,0.2983285084902258,0.2983285084902258,339.467141925856,3055.204277332704,0.1396843135690518,0.1396843135690518),
  @("ECs","Npnt","Itgb1","FAPs",3.0,1.0,2.128539,6.385617,0.4682231486221785,0.4682231486221785,3.0,1.0,172.558497,517.675491,0.3227862111630279,0.3227862111630279,367.297490645883,3305.677415812946,0.1511359761225763,0.1511359761225763),
  @("ECs","Npnt","Itgb1","Inflammatory-Mac",3.0,1.0,2.128539,6.385617,0.4682231486221785,0.4682231486221785,3.0,1.0,74.38770566666666,223.163117,0.1391489036280481,0.1391489036280482,158.337132632021,1425.034193688189,0.06515273778404877,0.06515273778404879),
  @("ECs","Npnt","Itgb1","MuSCs",3.0,1.0,2.128539,6.385617,0.4682231486221785,0.4682231486221785,3.0,1.0,58.41461433333333,175.243843,0.1092697975759847,0.1092697975759848,124.337784778459,1119.040063006131,0.05116264867033566,0.05116264867033567),
  @("ECs","Npnt","Itgb1","Resolving-Mac",3.0,1.0,2.128539,6.385617,0.4682231486221785,0.4682231486221785,3.0,1.0,69.746216,209.238648,0.1304665791427133,0.1304665791427133,148.457540858424,1336.117867725816,0.06108747247616586,0.06108747247616587),
  @("FAPs","Npnt","Itgb1","ECs",3.0,1.0,0.6061233333333333,1.81837,0.133331348679401,0.133331348679401,3.0,1.0,159.4836373333333,478.450912,0.2983285084902258,0.2983285084902258,96.66675387260446,870.00078485344,0.03977654238651594,0.03977654238651594),
  @("FAPs","Npnt","Itgb1","FAPs",3.0,1.0,0.6061233333333333,1.81837,0.133331348679401,0.133331348679401,3.0,1.0,172.558497,517.675491,0.3227862111630279,0.3227862111630279,104.59173139663,941.3255825696699,0.04303752086948043,0.04303752086948043),
  @("FAPs","Npnt","Itgb1","Inflammatory-Mac",3.0,1.0,0.6061233333333333,1.81837,0.133331348679401,0.133331348679401,3.0,1.0,74.38770566666666,223.163117,0.1391489036280481,0.1391489036280482,45.08812411769889,405.79311705929,0.01855291098798765,0.01855291098798766),
  @("FAPs","Npnt","Itgb1","MuSCs",3.0,1.0,0.6061233333333333,1.81837,0.133331348679401,0.133331348679401,3.0,1.0,58.41461433333333,175.243843,0.1092697975759847,0.1092697975759848,35.40646075510112,318.65814679591,0.01456908948073119,0.01456908948073119),
  @("FAPs","Npnt","Itgb1","Resolving-Mac",3.0,1.0,0.6061233333333333,1.81837,0.133331348679401,0.133331348679401,3.0,1.0,69.746216,209.238648,0.1304665791427133,0.1304665791427133,42.27480892930667,380.47328036376,0.01739528495468578,0.01739528495468578),
  @("MuSCs","Npnt","Itgb1","ECs",3.0,1.0,1.392473666666667,4.177421,0.306307943890216,0.306307943890216,3.0,1.0,159.4836373333333,478.450912,0.2983285084902258,0.2983285084902258,222.0767652508836,1998.690887257952,0.0913803920394759,0.0913803920394759),
  @("MuSCs","Npnt","Itgb1","FAPs",3.0,1.0,1.392473666666667,4.177421,0.306307943890216,0.306307943890216,3.0,1.0,172.558497,517.675491,0.3227862111630279,0.3227862111630279,240.283163032079,2162.548467288711,0.09887198065746017,0.09887198065746017),
  @("MuSCs","Npnt","Itgb1","Inflammatory-Mac",3.0,1.0,1.392473666666667,4.177421,0.306307943890216,0.306307943890216,3.0,1.0,74.38770566666666,223.163117,0.1391489036280481,0.1391489036280482,103.5829212645841,932.2462913812569,0.04262241456488524,0.04262241456488525),
  @("MuSCs","Npnt","Itgb1","MuSCs",3.0,1.0,1.392473666666667,4.177421,0.306307943890216,0.306307943890216,3.0,1.0,58.41461433333333,175.243843,0.1092697975759847,0.1092697975759848,81.34081220765589,732.067309868903,0.0334702070248,0.0334702070248),
  @("MuSCs","Npnt","Itgb1","Resolving-Mac",3.0,1.0,1.392473666666667,4.177421,0.306307943890216,0.306307943890216,3.0,1.0,69.746216,209.238648,0.1304665791427133,0.1304665791427133,97.11976912964533,874.077922166808,0.03996294960359465,0.03996294960359466),
  @("Resolving-Mac","Npnt","Itgb1","ECs",2.0,0.6666666666666666,0.4188566666666667,1.25657,0.09213755880820455,0.09213755880820455,3.0,1.0,159.4836373333333,478.450912,0.2983285084902258,0.2983285084902258,66.80078472131557,601.20706249184,0.02748726049518213,0.02748726049518213),
  @("Resolving-Mac","Npnt","Itgb1","FAPs",2.0,0.6666666666666666,0.4188566666666667,1.25657,0.09213755880820455,0.09213755880820455,3.0,1.0,172.558497,517.675491,0.3227862111630279,0.3227862111630279,72.27727685843,650.4954917258699,0.02974073351351102,0.02974073351351102),
  @("Resolving-Mac","Npnt","Itgb1","Inflammatory-Mac",2.0,0.6666666666666666,0.4188566666666667,1.25657,0.09213755880820455,0.09213755880820455,3.0,1.0,74.38770566666666,223.163117,0.1391489036280481,0.1391489036280482,31.15778643652111,280.42007792869,0.01282084029112647,0.01282084029112648),
  @("Resolving-Mac","Npnt","Itgb1","MuSCs",2.0,0.6666666666666666,0.4188566666666667,1.25657,0.09213755880820455,0.09213755880820455,3.0,1.0,58.41461433333333,175.243843,0.1092697975759847,0.1092697975759848,24.46735064427889,220.20615579851,0.0100678524001179,0.0100678524001179),
  @("Resolving-Mac","Npnt","Itgb1","Resolving-Mac",2.0,0.6666666666666666,0.4188566666666667,1.25657,0.09213755880820455,0.09213755880820455,3.0,1.0,69.746216,209.238648,0.1304665791427133,0.1304665791427133,29.21366754637333,262.92300791736,0.01202087210826702,0.01202087210826702),
)

$numRows = $rows.Count
$numCols = 20
$data = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
  for ($j = 0; $j -lt $numCols; $j++) {
    $data[$i,$j] = $rows[$i][$j]
  }
}

$ws.Range("A2:T21").Value = $data
